$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new Delivery Station / Region rows right after the last
# existing row of data (row 119 -> new rows 120 and 121).
$ws.Range("A120").Value = "ALK DS05"
$ws.Range("B120").Value = "East Region"

$ws.Range("A121").Value = "TBJ DS01"
$ws.Range("B121").Value = "North Region"

# Scroll the view down and select the newly added last cell, matching
# the saved sheetView state in the target workbook.
$ws.Activate()
$ws.Range("A121").Select()
$excel.ActiveWindow.ScrollRow = 99
